# Journal de bord - Anthony: mark up technical terms with spell-check
# proofErr tags, and append the November/December entries.

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Insert-ParaXml($range, [string]$bodyXml) {
    $range.InsertXML($pkgOpen + $bodyXml + $pkgClose)
}

# --- 1. "26" entry: split "métadatas" out with proofErr marks ---
$p = $d.Paragraphs.Item(3)
$body = '<w:p w:rsidR="003760D9" w:rsidRDefault="003760D9" w:rsidP="003760D9"><w:r><w:rPr><w:b/></w:rPr><w:t>26</w:t></w:r>' + `
  '<w:r><w:tab/></w:r>' + `
  '<w:r><w:tab/><w:t xml:space="preserve">- Test de diverse APIs pour lire les </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>m' + [char]0xE9 + 'tadatas</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> des MP3</w:t></w:r>' + `
  '</w:p>'
Insert-ParaXml $p.Range $body

# --- 2. Spotify entry: split "Spotify" out with proofErr marks ---
$p = $d.Paragraphs.Item(4)
$body = '<w:p w:rsidR="003760D9" w:rsidRDefault="003760D9" w:rsidP="003760D9"><w:r><w:tab/></w:r>' + `
  '<w:r><w:tab/><w:t xml:space="preserve">- Test de l' + [char]0x2019 + 'API </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>Spotify</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>'
Insert-ParaXml $p.Range $body

# --- 3. "Creation du SyncManager" entry: split "SyncManager" out ---
$p = $d.Paragraphs.Item(10)
$body = '<w:p w:rsidR="00026BF1" w:rsidRDefault="00026BF1"><w:r><w:tab/></w:r>' + `
  '<w:r><w:tab/><w:t xml:space="preserve">- Cr' + [char]0xE9 + 'ation du </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>SyncManager</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>'
Insert-ParaXml $p.Range $body

# --- 4. "Javadoc pour la partie TCP" entry: split "Javadoc" out ---
$p = $d.Paragraphs.Item(12)
$body = '<w:p w:rsidR="00A667F8" w:rsidRDefault="00BE081D"><w:r><w:tab/></w:r>' + `
  '<w:r><w:tab/><w:t xml:space="preserve">- </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>Javadoc</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> pour la partie TCP</w:t></w:r>' + `
  '</w:p>'
Insert-ParaXml $p.Range $body

# --- 5. "Tests et corrections de bug dans SyncManager": split "SyncManager" out ---
$p = $d.Paragraphs.Item(13)
$body = '<w:p w:rsidR="004B3E86" w:rsidRDefault="004B3E86"><w:r><w:rPr><w:b/></w:rPr><w:t>16</w:t></w:r>' + `
  '<w:r><w:tab/></w:r>' + `
  '<w:r><w:tab/><w:t xml:space="preserve">- Tests et corrections de bug dans </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>SyncManager</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '</w:p>'
Insert-ParaXml $p.Range $body

# --- 6. "Slide API et partie TCP" entry: drop the bookmark from here and
#        append the new November/December journal entries, ending with a
#        blank line and a paragraph that now carries the _GoBack bookmark.
$p = $d.Paragraphs.Item(14)
$body = '<w:p w:rsidR="00646FC4" w:rsidRPr="004B3E86" w:rsidRDefault="00646FC4"><w:r><w:tab/></w:r><w:r><w:tab/><w:t>- Slide API et partie TCP</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>23</w:t></w:r><w:r><w:tab/></w:r>' + `
    '<w:r><w:tab/><w:t>- Impl' + [char]0xE9 + 'mentation d' + [char]0x2019 + 'un contr' + [char]0xF4 + 'leur central</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>30</w:t></w:r><w:r><w:tab/></w:r>' + `
    '<w:r><w:tab/><w:t>- Documentation</w:t></w:r></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Titre1"/></w:pPr><w:r><w:t>D' + [char]0xE9 + 'cembre 2015</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>7</w:t></w:r><w:r><w:tab/></w:r>' + `
    '<w:r><w:tab/><w:t>- Mise en place de la synchronisation du lecteur audio</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>14</w:t></w:r><w:r><w:tab/></w:r>' + `
    '<w:r><w:tab/><w:t>- Correction de bugs lors de la synchronisation</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>- Documentation</w:t></w:r></w:p>' + `
  '<w:p/>' + `
  '<w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Insert-ParaXml $p.Range $body
